$wb = $excel.ActiveWorkbook

# Update F3, F5, F8 on both the "展览" and "全部类型" sheets
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1831
    $ws.Range("F5").Value = 1126
    $ws.Range("F8").Value = 5929
}
